# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型"
# sheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 224
    4  = 12999
    5  = 1341
    6  = 215
    7  = 38
    9  = 165
    17 = 411
    18 = 5533
    20 = 54
    24 = 139
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
